$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update legmaxROM / meanEMG values (row 1: counts, rows 2-3: CON/STR data) in columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 0.4763908301733073
$ws.Range("C2").Value = 6.3118189179278641
$ws.Range("D2").Value = 7.9502199932933806
$ws.Range("E2").Value = 12.184128701107621

$ws.Range("B3").Value = -5.9463987840165569
$ws.Range("C3").Value = 5.2817447298443208
$ws.Range("D3").Value = 15.634886854866414
$ws.Range("E3").Value = 8.0830055985159763

# Update the selection to reflect the now-smaller range of interest
$ws.Range("B1:E3").Select()
